$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in columns A, B, E, F, G, Q, R between row 4 and row 5.
$cols = @("A", "B", "E", "F", "G", "Q", "R")

foreach ($col in $cols) {
    $addr4 = "${col}4"
    $addr5 = "${col}5"
    $val4 = $ws.Range($addr4).Value2
    $val5 = $ws.Range($addr5).Value2
    $ws.Range($addr4).Value = $val5
    $ws.Range($addr5).Value = $val4
}
